$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rearrange columns A:E (Campaign, Ad group, Keyword, Match type, New Bid)
# into the new order (Keyword, New Bid, Campaign, Ad group, Match type)
# for the header row and the two data rows, preserving values/types.

# Row 1 (headers)
$ws.Range("A1").Value = "Keyword"
$ws.Range("B1").Value = "New Bid"
$ws.Range("C1").Value = "Campaign"
$ws.Range("D1").Value = "Ad group"
$ws.Range("E1").Value = "Match type"

# Row 2
$ws.Range("A2").Value = "new homes for sale spring tx"
$ws.Range("B2").Value = 0.4
$ws.Range("C2").Value = "Houston_TX>279>SB>City_MSM102"
$ws.Range("D2").Value = "nonmarketname>newhometerms>Spring_TX>279"
$ws.Range("E2").Value = "Broad"

# Row 3
$ws.Range("A3").Value = "rosharon tx new homes"
$ws.Range("B3").Value = 0.18
$ws.Range("C3").Value = "Houston_TX>279>SB>City_MSM102"
$ws.Range("D3").Value = "nonmarketname>newhometerms>Rosharon_TX>279"
$ws.Range("E3").Value = "Broad"

# Update the selection to match the committed sheet view state.
$ws.Range("A1:T1").Select()
